# FuelPrices update: insert a new latest-date row at the top of the data
# table (row 2), shifting all existing data rows down by one, and append
# the row that fell off the bottom of the 28-row window as a fresh blank
# row (price not yet published) with the oldest tracked date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (above the previous most-recent entry),
# pushing every existing row down by one. Excel shifts cell content
# (values + existing formatting) down together, so the whole history
# -- including already-blank cells -- stays intact.
$ws.Rows("2:2").Insert()

# Populate the freshly inserted row with the newest reading.
$ws.Cells.Item(2, 1).Value2 = 45769
$ws.Cells.Item(2, 2).Value2 = 748.228
$ws.Cells.Item(2, 3).Value2 = 753.256

# The Insert() operation stamped row 2 with a copied/derived style that
# doesn't match the rest of the date column. Fix formatting by copying
# it from the row directly below (which holds the correctly-styled data
# that used to be row 2, i.e. date-formatted column A / plain B & C).
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Inserted new FuelPrices row for 2025-04-22 (serial 45769)."
